$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("fixed")

# --- Add column D ("regex_capture_full") values to the "fixed" sheet ---
$ws2.Range("D1").Value = 'regex_capture_full'
$ws2.Range("D2").Value = 'antarctica'
$ws2.Range("D6").Value = 'bouvet\\w*( island\\w*)?'
$ws2.Range("D7").Value = 'british.?indian.?ocean( territor\\w*)?'
$ws2.Range("D9").Value = 'christmas island\\w*'
$ws2.Range("D10").Value = 'cocos island\\w*|keeling'
$ws2.Range("D11").Value = 'czechoslovak\\w*'
$ws2.Range("D12").Value = 'french.?southern( and |[ &\\/)]*)antarct\\w*(lands)?'
$ws2.Range("D13").Value = 'german.?democratic.?rep|democratic.?rep.*germany|east.germany|germany[-\\ (]east|\\bgdr\\b'
$ws2.Range("D16").Value = 'heard (island)?( and |[ &\\/)]*)mcdonald( island\\w*)'
$ws2.Range("D23").Value = 'netherlands.antil\\w*|dutch.antil\\w*'
$ws2.Range("D33").Value = 'serbia( and |[ &\\/)]*)montenegro'
$ws2.Range("D34").Value = 'somaliland\\w*'
$ws2.Range("D35").Value = 'south.?georgia|sandwich island\\w*|south.?georgia( and |[ &\\/)]*)(the )?sandwich( island\\w*)?'
$ws2.Range("D37").Value = 'tibet\\w*'
$ws2.Range("D42").Value = '(united states |u\\.?s\\.? )?minor.?outlying.?is\\w*'
$ws2.Range("D47").Value = '^(?=.*peo).*yemen|^(?!.*rep)(?=.*dem).*yemen|^(?=.*south).*yemen|^(?=.*aden).*yemen|^(?=.*\\bp\.?d\.?r).*yemen'
$ws2.Range("D48").Value = '(?<!former )yugoslav\\w*'

# --- Update column B regex values that were revised ---
$ws2.Range("B13").Value = 'german.?democratic.?rep|democratic.?rep.*germany|east.germany|germany[-\\ (]east|\\bgdr\\b'
$ws2.Range("B16").Value = 'heard (island)?( and |[ &\\/)]*)mcdonald'
$ws2.Range("B48").Value = '(?<!former )yugoslav'

# --- Restore the selection / active-cell view state captured in the workbook ---
$ws1.Activate()
$ws1.Range("C30").Select() | Out-Null
$ws2.Activate()
$ws2.Range("D6").Select() | Out-Null
